$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 424.45456
$ws.Range("I4").Value = 129.71428
$ws.Range("K4").Value = 129.71428
$ws.Range("M4").Value = -15.71428

$ws.Range("H40").Value = 6990.1
$ws.Range("J40").Value = 8050.25
$ws.Range("L40").Value = 8050.25
$ws.Range("N40").Value = -8400.25

$ws.Range("H43").Value = 1964.9
$ws.Range("I43").Value = 2141.6667
$ws.Range("K43").Value = 2141.6667
$ws.Range("M43").Value = -2072.6667

$ws.Range("H64").Value = 6168.6113
$ws.Range("I64").Value = 5207.3
$ws.Range("J64").Value = 7370.25
$ws.Range("K64").Value = 5207.3
$ws.Range("L64").Value = 7370.25
$ws.Range("M64").Value = -4959.3
$ws.Range("N64").Value = -7866.25

$ws.Range("H67").Value = 6168.6113
$ws.Range("I67").Value = 5207.3
$ws.Range("J67").Value = 7370.25
$ws.Range("K67").Value = 5207.3
$ws.Range("L67").Value = 7370.25
$ws.Range("M67").Value = -4349.3
$ws.Range("N67").Value = -9086.25

$ws.Range("H100").Value = 4954.8667
$ws.Range("I100").Value = 3390
$ws.Range("J100").Value = 8084.6
$ws.Range("K100").Value = 3390
$ws.Range("L100").Value = 8084.6
$ws.Range("M100").Value = -2849
$ws.Range("N100").Value = -9166.6

$ws.Range("H112").Value = 3448.9565
$ws.Range("J112").Value = 3158.476
$ws.Range("L112").Value = 9475.428
$ws.Range("N112").Value = -11691.428

$ws.Range("H116").Value = 6747.5
$ws.Range("I116").Value = 9000
$ws.Range("K116").Value = 9000
$ws.Range("M116").Value = -5558

$ws.Range("H135").Value = 1532.2916
$ws.Range("I135").Value = 1532.2916
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13790.6244
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11255.6244
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 2322.0952
$ws.Range("I138").Value = 1909.5
$ws.Range("J138").Value = 2576
$ws.Range("K138").Value = 5728.5
$ws.Range("L138").Value = 7728
$ws.Range("M138").Value = -588.5
$ws.Range("N138").Value = -18008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4689.75
$ws.Range("J5").Value = 12999.4
$ws.Range("L5").Value = 12999.4
$ws.Range("N5").Value = -13223.4

$ws.Range("H32").Value = 6427.4
$ws.Range("I32").Value = 4314.4185
$ws.Range("K32").Value = 4314.4185
$ws.Range("M32").Value = -4027.4185

$ws.Range("H61").Value = 2361.5652
$ws.Range("I61").Value = 2014.8649
$ws.Range("K61").Value = 2014.8649
$ws.Range("M61").Value = -1802.8649

$ws.Range("H74").Value = 1522.1892
$ws.Range("I74").Value = 1499.1562
$ws.Range("J74").Value = 1669.6
$ws.Range("K74").Value = 1499.1562
$ws.Range("L74").Value = 1669.6
$ws.Range("M74").Value = -625.1561999999999
$ws.Range("N74").Value = -3417.6

$ws.Range("H77").Value = 1522.1892
$ws.Range("I77").Value = 1499.1562
$ws.Range("J77").Value = 1669.6
$ws.Range("K77").Value = 7495.780999999999
$ws.Range("L77").Value = 8348
$ws.Range("M77").Value = -3127.780999999999
$ws.Range("N77").Value = -17084

$ws.Range("H109").Value = 352782.75
$ws.Range("J109").Value = 352782.75
$ws.Range("L109").Value = 352782.75
$ws.Range("N109").Value = -355556.75

$ws.Range("H130").Value = 35214.5
$ws.Range("J130").Value = 35214.5
$ws.Range("L130").Value = 35214.5
$ws.Range("N130").Value = -45254.5

$ws.Range("H136").Value = 2361.5652
$ws.Range("I136").Value = 2014.8649
$ws.Range("K136").Value = 6044.5947
$ws.Range("M136").Value = -3494.5947

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4689.75
$ws.Range("J4").Value = 12999.4
$ws.Range("L4").Value = 12999.4
$ws.Range("N4").Value = -13229.4

$ws.Range("H22").Value = 327.57144
$ws.Range("J22").Value = 297
$ws.Range("L22").Value = 297
$ws.Range("N22").Value = -643

$ws.Range("H92").Value = 12401
$ws.Range("J92").Value = 12401
$ws.Range("L92").Value = 12401
$ws.Range("N92").Value = -17393

$ws.Range("H122").Value = 99999
$ws.Range("J122").Value = 99999
$ws.Range("L122").Value = 99999
$ws.Range("N122").Value = -109799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2124.561
$ws.Range("I31").Value = 1056.04
$ws.Range("J31").Value = 3794.125
$ws.Range("K31").Value = 1056.04
$ws.Range("L31").Value = 3794.125
$ws.Range("M31").Value = -761.04
$ws.Range("N31").Value = -4384.125

$ws.Range("H34").Value = 2124.561
$ws.Range("I34").Value = 1056.04
$ws.Range("J34").Value = 3794.125
$ws.Range("K34").Value = 1056.04
$ws.Range("L34").Value = 3794.125
$ws.Range("M34").Value = -854.04
$ws.Range("N34").Value = -4198.125

$ws.Range("H52").Value = 125637.5
$ws.Range("J52").Value = 125637.5
$ws.Range("L52").Value = 125637.5
$ws.Range("N52").Value = -126225.5

$ws.Range("H107").Value = 733.13794
$ws.Range("I107").Value = 420.73685
$ws.Range("K107").Value = 420.73685
$ws.Range("M107").Value = 1499.26315

$ws.Range("H132").Value = 2656.5
$ws.Range("I132").Value = 2399.3076
$ws.Range("K132").Value = 7197.9228
$ws.Range("M132").Value = -4667.9228

$ws.Range("H134").Value = 1941.2826
$ws.Range("I134").Value = 1494.4412
$ws.Range("K134").Value = 4483.3236
$ws.Range("M134").Value = -1948.3236

$ws.Range("H139").Value = 59678
$ws.Range("J139").Value = 59678
$ws.Range("L139").Value = 59678
$ws.Range("N139").Value = -69958

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H131").Value = 3176430
$ws.Range("I131").Value = 10001986
$ws.Range("K131").Value = 30005958
$ws.Range("M131").Value = -30000918

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 19199.834
$ws.Range("J3").Value = 39999.5
$ws.Range("L3").Value = 39999.5
$ws.Range("N3").Value = -40231.5

$ws.Range("H19").Value = 50000
$ws.Range("I19").Value = 50000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 50000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -49712
$ws.Range("N19").ClearContents()

$ws.Range("H97").Value = 412.34784
$ws.Range("I97").Value = 403.625
$ws.Range("J97").Value = 432.2857
$ws.Range("K97").Value = 403.625
$ws.Range("L97").Value = 432.2857
$ws.Range("M97").Value = 92.375
$ws.Range("N97").Value = -1424.2857

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H113").Value = 3250
$ws.Range("I113").Value = 3289.6428
$ws.Range("J113").Value = 3188.3333
$ws.Range("K113").Value = 3289.6428
$ws.Range("L113").Value = 3188.3333
$ws.Range("M113").Value = -1119.6428
$ws.Range("N113").Value = -7528.3333

$ws.Range("H132").Value = 3520.862
$ws.Range("I132").Value = 3362.9023
$ws.Range("J132").Value = 3901.8235
$ws.Range("K132").Value = 10088.7069
$ws.Range("L132").Value = 11705.4705
$ws.Range("M132").Value = -7558.706900000001
$ws.Range("N132").Value = -16765.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2993.923
$ws.Range("I136").Value = 2186.4
$ws.Range("K136").Value = 6559.200000000001
$ws.Range("M136").Value = -4009.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 45570.82
$ws.Range("I62").Value = 61027.1
$ws.Range("J62").Value = 6930.125
$ws.Range("K62").Value = 61027.1
$ws.Range("L62").Value = 6930.125
$ws.Range("M62").Value = -60403.1
$ws.Range("N62").Value = -8178.125

$ws.Range("H65").Value = 45570.82
$ws.Range("I65").Value = 61027.1
$ws.Range("J65").Value = 6930.125
$ws.Range("K65").Value = 305135.5
$ws.Range("L65").Value = 34650.625
$ws.Range("M65").Value = -302015.5
$ws.Range("N65").Value = -40890.625

$ws.Range("H100").Value = 468
$ws.Range("J100").Value = 518.75
$ws.Range("L100").Value = 1037.5
$ws.Range("N100").Value = -2119.5

$ws.Range("H135").Value = 83571.336
$ws.Range("J135").Value = 83571.336
$ws.Range("L135").Value = 83571.336
$ws.Range("N135").Value = -93711.336
